$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme
try {
  $tcs.Load("Office")
  Write-Output "load ok"
} catch {
  Write-Output "load err: $_"
}
try {
  $tcs.Save("Office")
  Write-Output "save ok"
} catch {
  Write-Output "save err: $_"
}
